$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Childcare SP Account name",
    "Little Cubs",
    "Keith One",
    "Little Cubs",
    "Keith One",
    "Childcare SP Account name",
    "Little Cubs",
    "ates"
)

$startRow = 56
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
